{"js": "// Remove the \"CreateDeck (abstract)\" line from the Deck Class (Abstract)\n// table's \"Methods\" cell, leaving only \"Shuffle (void)\" behind.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"CreateDeck\") !== -1 && p.text.indexOf(\"abstract\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'CreateDeck (abstract)' paragraph.\");\n}\n\ntarget.delete();\nawait context.sync();\n", "ps1": "# Remove the \"CreateDeck (abstract)\" line from the Deck Class (Abstract)\n# table's \"Methods\" cell, leaving only \"Shuffle (void)\" behind.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*CreateDeck*\" -and $t -like \"*abstract*\") {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'CreateDeck (abstract)' paragraph.\"\n}\n\n$target.Range.Delete()\n"}
